$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Playbook1"

# Header row (row 1) - write in the order that reproduces the target
# shared-string table ordering.
$ws.Range("B1").Value = "IP"
$ws.Range("C1").Value = "Metering"
$ws.Range("E1").Value = "VLAN"
$ws.Range("B2").Value = "10.0.252.11"
$ws.Range("C2").Value = "Gi1/0/21 - 23"
$ws.Range("F1").Value = "CoreIP"
$ws.Range("F2").Value = "10.0.252.1"
$ws.Range("G1").Value = "CoreTrunk"
$ws.Range("D1").Value = "IDFTrunk"
$ws.Range("A1").Value = "Floor"

# Row 2 numeric cell
$ws.Range("E2").Value = 16

# Clear old D2 value (170) since D column data moved
$ws.Range("D2").Value = ""
